$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$symbol = ":ECBASSETSW"

# --- Fill in the missing symbol (column B) for existing rows 1327 and 1328 ---
$ws.Range("B1327").Value = $symbol
$ws.Range("B1328").Value = $symbol

# --- Apply the same date style used by A1328 to the new A-column cells ---
$ws.Range("A1328").Copy()
$ws.Range("A1329:A1345").PasteSpecial(-4122)  # xlPasteFormats

# --- New row data (date serial, symbol flag, value) ---
$rowData = @(
    @{ Row = 1329; Date = 45261; HasSymbol = $true;  Value = 7002047000000 },
    @{ Row = 1330; Date = 45268; HasSymbol = $true;  Value = 6993472000000 },
    @{ Row = 1331; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1332; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1333; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1334; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1335; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1336; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1337; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1338; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1339; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1340; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1341; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1342; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1343; Date = 45271; HasSymbol = $true;  Value = 6987406000000 },
    @{ Row = 1344; Date = 45278; HasSymbol = $true;  Value = 6899165000000 },
    @{ Row = 1345; Date = 45278; HasSymbol = $false; Value = 6899165000000 }
)

foreach ($r in $rowData) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Date
    if ($r.HasSymbol) {
        $ws.Range("B$row").Value = $symbol
    }
    $ws.Range("C$row").Value = $r.Value
    $ws.Range("D$row").Value = $r.Value
    $ws.Range("E$row").Value = $r.Value
    $ws.Range("F$row").Value = $r.Value
    $ws.Range("G$row").Value = 0
}
